$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest Handoff Date" (column D) rows 7-16 -> overall latest handoff timestamp
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D7:D16").Value = "2016-05-18 07:05:36"

# zh-cn sheet: "Latest Handoff Datetime" (column E) rows 7-16 -> new handoff timestamp for zh-cn
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E7:E16").Value = "2016-03-18 07:05:32"

# de-de sheet: "Latest Handoff Datetime" (column E) rows 7-16 -> new handoff timestamp for de-de
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E7:E16").Value = "2016-03-18 07:05:36"
